$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.349.63'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.712.62'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  +0.54%  '
$ws.Range("D5").Value = '''224.42'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").Value = '''0.5263'
$ws.Range("E6").Value = '  -0.41%  '
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").Value = '''0.06628'
$ws.Range("E8").Value = '  +1.38%  '
$ws.Range("D9").Value = '''0.2638'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("D10").Value = '''20.69'
$ws.Range("E10").Value = '  -1.09%  '
$ws.Range("D11").Value = '''0.07746'
$ws.Range("E11").Value = '  +1.53%  '
$ws.Range("D12").Value = '''4.452'
$ws.Range("E12").Value = '  -2.46%  '
$ws.Range("D13").Value = '1.950.36'
$ws.Range("E13").Value = '  +0.31%  '
$ws.Range("D14").Value = '1.710.37'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '''0.5764'
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("D16").Value = '0.0₅8167'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").Value = '''67.62'
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = '27.374.96'
$ws.Range("D19").Value = '''218.82'
$ws.Range("E19").Value = '  +1.51%  '
$ws.Range("E20").Value = '  +0.53%  '
$ws.Range("D21").Value = '''4.635'
$ws.Range("E21").Value = '  -0.90%  '
$ws.Range("D22").Value = '''10.41'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = '''6.018'
$ws.Range("E23").Value = '  +1.09%  '
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = '''145.39'
$ws.Range("E25").Value = '  +2.01%  '
$ws.Range("D26").Value = '''1.726'
$ws.Range("E26").Value = '  -1.49%  '
$ws.Range("D27").Value = '''0.1200'
$ws.Range("E27").Value = '  -1.30%  '
$ws.Range("D28").Value = '''7.197'
$ws.Range("E28").Value = '  -0.89%  '
$ws.Range("D29").Value = '''16.13'
$ws.Range("E29").Value = '  -0.98%  '
$ws.Range("D30").Value = '''0.05305'
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D32").Value = '''3.470'
$ws.Range("E32").Value = '  -0.38%  '
$ws.Range("D33").Value = '''3.349'
$ws.Range("E33").Value = '  -1.83%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = '''2.834'
$ws.Range("E35").Value = '  -1.21%  '
$ws.Range("D36").Value = '''0.9484'
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").Value = '''2.402'
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").Value = '''0.5863'
$ws.Range("E38").Value = '  +0.04%  '
$ws.Range("D39").Value = '1.186.71'
$ws.Range("E39").Value = '  +14.06%  '
$ws.Range("D40").Value = '''0.01648'
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").Value = '''5.785'
$ws.Range("E41").Value = '  -1.30%  '
$ws.Range("E42").Value = '  +0.49%  '
$ws.Range("D43").Value = '''0.8394'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '''101.07'
$ws.Range("E44").Value = '  +0.08%  '
$ws.Range("D45").Value = '1.857.48'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("E46").Value = '  +3.13%  '
$ws.Range("D47").Value = '''57.39'
$ws.Range("E47").Value = '  -0.99%  '
$ws.Range("D48").Value = '''0.4562'
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("E49").Value = '  +0.40%  '
$ws.Range("D50").Value = '''8.137'
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("D51").Value = '''0.05236'
$ws.Range("E51").Value = '  -0.09%  '
